$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Agen" column (column I), which shifts everything after it left by one.
$ws.Columns("I").Delete()

# Update the active selection to I1 (matches post-edit cursor location in the source file).
$ws.Range("I1").Select()
